$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (4-10) of data get reshuffled (a single 7-cycle permutation of the
# record rows), which changes the visible values in columns A,B,D,E,F,G,H,Q,R,AC
# while all other columns remain identical across the group (location/date/etc).
# We simply rewrite the affected columns with their final values.

$rows = @{
  4  = @{ A = 86419304; B = 90653;  D = "LC"; E = 4364;   F = "Dropptaggsvamp";   G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."; Q = 440783.8270494898; R = 6707144.091754919; AC = $null }
  5  = @{ A = 86419294; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";          G = "Alectoria sarmentosa";  H = "(Ach.) Ach.";        Q = 440686.1394479795; R = 6707147.171128325; AC = "På tall" }
  6  = @{ A = 86419305; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";          G = "Alectoria sarmentosa";  H = "(Ach.) Ach.";        Q = 440606.8734944779; R = 6707280.052989913; AC = "På flera tallar" }
  7  = @{ A = 86419296; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";          G = "Alectoria sarmentosa";  H = "(Ach.) Ach.";        Q = 440814.1817916233; R = 6707128.810180089; AC = "På tall" }
  8  = @{ A = 86419290; B = 8377;   D = "LC"; E = 106545; F = "Mindre märgborre"; G = "Tomicus minor";         H = "(Hartig, 1834)";     Q = 440814.1656648018; R = 6707127.824798071; AC = $null }
  9  = @{ A = 86419313; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";          G = "Alectoria sarmentosa";  H = "(Ach.) Ach.";        Q = 440607.1726549119; R = 6707238.159541409; AC = "rikligt på flera tallar" }
  10 = @{ A = 86419293; B = 77506;  D = "NT"; E = 6425;   F = "Garnlav";          G = "Alectoria sarmentosa";  H = "(Ach.) Ach.";        Q = 440687.1425972193; R = 6707148.140317255; AC = $null }
}

foreach ($r in $rows.Keys) {
  $row = $rows[$r]
  $ws.Range("A$r").Value = $row.A
  $ws.Range("B$r").Value = $row.B
  $ws.Range("D$r").Value = $row.D
  $ws.Range("E$r").Value = $row.E
  $ws.Range("F$r").Value = $row.F
  $ws.Range("G$r").Value = $row.G
  $ws.Range("H$r").Value = $row.H
  $ws.Range("Q$r").Value = $row.Q
  $ws.Range("R$r").Value = $row.R
  if ($row.AC -eq $null) {
    $ws.Range("AC$r").ClearContents()
  } else {
    $ws.Range("AC$r").Value = $row.AC
  }
}
